$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.6109724624622288
$ws.Range("C2").Value = 0.05997153799206956
$ws.Range("D2").Value = 0.05856964968633349
$ws.Range("F2").Value = 1.349328526400285
$ws.Range("G2").Value = 0.002489829265264314
$ws.Range("K2").Value = 0.5938401813831149
$ws.Range("N2").Value = 1.963158319943183
$ws.Range("B3").Value = 0.568226392564668
$ws.Range("C3").Value = 0.05307900734466386
$ws.Range("D3").Value = 0.05824567024120952
$ws.Range("F3").Value = 1.342611000863641
$ws.Range("G3").Value = 0.002493334951751594
$ws.Range("K3").Value = 0.5468139539168817
$ws.Range("N3").Value = 1.982660104844271
$ws.Range("B4").Value = 0.5422921543181189
$ws.Range("C4").Value = 0.0488611280244271
$ws.Range("D4").Value = 0.05804512657321226
$ws.Range("F4").Value = 1.339284429731222
$ws.Range("G4").Value = 0.002495601496713444
$ws.Range("K4").Value = 0.5182233303110877
$ws.Range("N4").Value = 1.995314520940575
$ws.Range("B5").Value = 0.5318022431780207
$ws.Range("C5").Value = 0.04714583641532499
$ws.Range("D5").Value = 0.05796300782772335
$ws.Range("F5").Value = 1.338129170954666
$ws.Range("G5").Value = 0.002496553897392021
$ws.Range("K5").Value = 0.5066437329278415
$ws.Range("N5").Value = 2.000642163209513
$ws.Range("B6").Value = 0.5300651440662136
$ws.Range("C6").Value = 0.04686122572331897
$ws.Range("D6").Value = 0.05794934852337263
$ws.Range("F6").Value = 1.337949432168131
$ws.Range("G6").Value = 0.002496713782854378
$ws.Range("K6").Value = 0.5047252564274345
$ws.Range("N6").Value = 2.001537128541678
$ws.Range("B7").Value = 0.5421503658014046
$ws.Range("C7").Value = 0.04883798075672985
$ws.Range("D7").Value = 0.05804402067745329
$ws.Range("F7").Value = 1.339268038744152
$ws.Range("G7").Value = 0.002495614224425003
$ws.Range("K7").Value = 0.5180668750837185
$ws.Range("N7").Value = 1.99538567981698
$ws.Range("B8").Value = 0.5961688668097054
$ws.Range("C8").Value = 0.05759203205712993
$ws.Range("D8").Value = 0.05845828435034761
$ws.Range("F8").Value = 1.346846495239447
$ws.Range("G8").Value = 0.002491014410929775
$ws.Range("K8").Value = 0.5775666354073223
$ws.Range("N8").Value = 1.96974124014136
$ws.Range("B9").Value = 0.7045791812280413
$ws.Range("C9").Value = 0.074873448051477
$ws.Range("D9").Value = 0.05925734663540183
$ws.Range("F9").Value = 1.368056930979037
$ws.Range("G9").Value = 0.002482894958900579
$ws.Range("K9").Value = 0.6965037981217108
$ws.Range("N9").Value = 1.924857722016625
$ws.Range("B10").Value = 0.7857554829511173
$ws.Range("C10").Value = 0.08764464024568497
$ws.Range("D10").Value = 0.05983575069961589
$ws.Range("F10").Value = 1.387538509164017
$ws.Range("G10").Value = 0.002477473003175391
$ws.Range("K10").Value = 0.7852837813742326
$ws.Range("N10").Value = 1.8951866288483
$ws.Range("B11").Value = 0.8230201785676741
$ws.Range("C11").Value = 0.0934718400876875
$ws.Range("D11").Value = 0.06009688977066929
$ws.Range("F11").Value = 1.39725390089778
$ws.Range("G11").Value = 0.002475123190476948
$ws.Range("K11").Value = 0.8259803615659393
$ws.Range("N11").Value = 1.882408208735839
$ws.Range("B12").Value = 0.837179954690555
$ws.Range("C12").Value = 0.09568102677783941
$ws.Range("D12").Value = 0.06019548163601129
$ws.Range("F12").Value = 1.401055981142505
$ws.Range("G12").Value = 0.002474250061281305
$ws.Range("K12").Value = 0.8414358842540537
$ws.Range("N12").Value = 1.877672978657678
$ws.Range("B13").Value = 0.8341282396152678
$ws.Range("C13").Value = 0.09520512483686616
$ws.Range("D13").Value = 0.06017426144561355
$ws.Range("F13").Value = 1.400231655159971
$ws.Range("G13").Value = 0.002474437364105099
$ws.Range("K13").Value = 0.8381052765362824
$ws.Range("N13").Value = 1.878688179572862
$ws.Range("B14").Value = 0.8241841412039435
$ws.Range("C14").Value = 0.09365353985646152
$ws.Range("D14").Value = 0.06010500697504639
$ws.Range("F14").Value = 1.39756423151097
$ws.Range("G14").Value = 0.002475051023727518
$ws.Range("K14").Value = 0.8272510025463475
$ws.Range("N14").Value = 1.882016559020428
$ws.Range("B15").Value = 0.818099401495374
$ws.Range("C15").Value = 0.09270348266943529
$ws.Range("D15").Value = 0.06006254775921249
$ws.Range("F15").Value = 1.395946396881783
$ws.Range("G15").Value = 0.002475429079479841
$ws.Range("K15").Value = 0.8206082565113775
$ws.Range("N15").Value = 1.884068797951024
$ws.Range("B16").Value = 0.7833269288479698
$ws.Range("C16").Value = 0.08726417441110357
$ws.Range("D16").Value = 0.05981864384933644
$ws.Range("F16").Value = 1.38692078613019
$ws.Range("G16").Value = 0.002477628909515986
$ws.Range("K16").Value = 0.7826304058742437
$ws.Range("N16").Value = 1.896036222234216
$ws.Range("B17").Value = 0.7620814741288484
$ws.Range("C17").Value = 0.08393184581771607
$ws.Range("D17").Value = 0.05966850161784265
$ws.Range("F17").Value = 1.381602651413374
$ws.Range("G17").Value = 0.002479008255211995
$ws.Range("K17").Value = 0.7594116644403357
$ws.Range("N17").Value = 1.903562223645281
$ws.Range("B18").Value = 0.7498933906989294
$ws.Range("C18").Value = 0.08201682564384782
$ws.Range("D18").Value = 0.05958195807337319
$ws.Range("F18").Value = 1.378624077334209
$ws.Range("G18").Value = 0.002479812604201344
$ws.Range("K18").Value = 0.7460860473594835
$ws.Range("N18").Value = 1.907958649915216
$ws.Range("B19").Value = 0.7457721678229063
$ws.Range("C19").Value = 0.08136871460580153
$ws.Range("D19").Value = 0.0595526243651392
$ws.Range("F19").Value = 1.377629358911548
$ws.Range("G19").Value = 0.002480086832293918
$ws.Range("K19").Value = 0.7415792351565074
$ws.Range("N19").Value = 1.909458818329696
$ws.Range("B20").Value = 0.7643398043787215
$ws.Range("C20").Value = 0.08428640717761482
$ws.Range("D20").Value = 0.0596845038144842
$ws.Range("F20").Value = 1.382160464780341
$ws.Range("G20").Value = 0.002478860285176279
$ws.Range("K20").Value = 0.7618803180585303
$ws.Range("N20").Value = 1.902754062450896
$ws.Range("B21").Value = 0.8271036495103772
$ws.Range("C21").Value = 0.09410920843365034
$ws.Range("D21").Value = 0.06012535680774178
$ws.Range("F21").Value = 1.398344375197524
$ws.Range("G21").Value = 0.002474870324691897
$ws.Range("K21").Value = 0.830437955057107
$ws.Range("N21").Value = 1.881036117243646
$ws.Range("B22").Value = 0.8684058050031354
$ws.Range("C22").Value = 0.1005438710104158
$ws.Range("D22").Value = 0.06041175038199853
$ws.Range("F22").Value = 1.409639040118506
$ws.Range("G22").Value = 0.002472359920556076
$ws.Range("K22").Value = 0.8755043705588434
$ws.Range("N22").Value = 1.867446642385559
$ws.Range("B23").Value = 0.8463362453063041
$ws.Range("C23").Value = 0.09710819659215986
$ws.Range("D23").Value = 0.06025905859348413
$ws.Range("F23").Value = 1.403545078727291
$ws.Range("G23").Value = 0.002473690897506532
$ws.Range("K23").Value = 0.8514277770460978
$ws.Range("N23").Value = 1.87464420159808
$ws.Range("B24").Value = 0.7633187316250485
$ws.Range("C24").Value = 0.08412610756104755
$ws.Range("D24").Value = 0.05967726992040312
$ws.Range("F24").Value = 1.381908031634509
$ws.Range("G24").Value = 0.002478927147115464
$ws.Range("K24").Value = 0.7607641676630408
$ws.Range("N24").Value = 1.903119215244388
$ws.Range("B25").Value = 0.6749842302265279
$ws.Range("C25").Value = 0.07018561075551588
$ws.Range("D25").Value = 0.05904266201211428
$ws.Range("F25").Value = 1.361636378572101
$ws.Range("G25").Value = 0.002484995646420793
$ws.Range("K25").Value = 0.6640843074048064
$ws.Range("N25").Value = 2.000642163209513